$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.107.86"
Set-TextValue "E2" "  +0.22%  "

Set-TextValue "D3" "1.881.62"
Set-TextValue "E3" "  -1.00%  "

Set-TextValue "E4" "  +0.12%  "

Set-TextValue "D5" "313.47"
Set-TextValue "E5" "  +0.30%  "

Set-TextValue "E6" "  +0.08%  "

Set-TextValue "D7" "0.5077"
Set-TextValue "E7" "  +0.47%  "

Set-TextValue "D8" "0.3867"
Set-TextValue "E8" "  -1.43%  "

Set-TextValue "D9" "0.09002"
Set-TextValue "E9" "  -3.64%  "

Set-TextValue "D10" "1.125"
Set-TextValue "E10" "  -1.05%  "

Set-TextValue "E11" "  -0.28%  "

Set-TextValue "D12" "6.366"
Set-TextValue "E12" "  -0.08%  "

Set-TextValue "D13" "20.80"
Set-TextValue "E13" "  +0.15%  "

Set-TextValue "D14" "1.879.14"
Set-TextValue "E14" "  -0.36%  "

Set-TextValue "D15" "7.252"
Set-TextValue "E15" "  -0.79%  "

Set-TextValue "E16" "  +0.12%  "

Set-TextValue "D17" "0.00001110"
Set-TextValue "E17" "  -0.59%  "

Set-TextValue "D18" "91.38"
Set-TextValue "E18" "  -1.12%  "

Set-TextValue "D19" "0.06623"
Set-TextValue "E19" "  +0.70%  "

Set-TextValue "D20" "18.21"
Set-TextValue "E20" "  +2.21%  "

Set-TextValue "E21" "  +0.15%  "

Set-TextValue "E22" "  -1.39%  "

Set-TextValue "D23" "28.142.52"
Set-TextValue "E23" "  +0.12%  "

Set-TextValue "E24" "  +0.61%  "

Set-TextValue "D25" "2.272"
Set-TextValue "E25" "  -2.00%  "

Set-TextValue "B26" "LidoDAOToken"
Set-TextValue "C26" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D26" "2.551"
Set-TextValue "E26" "  -3.05%  "

Set-TextValue "D27" "2.089.56"
Set-TextValue "E27" "  -0.93%  "

Set-TextValue "B28" "EthereumClassic"
Set-TextValue "C28" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D28" "20.82"
Set-TextValue "E28" "  -0.25%  "

Set-TextValue "B29" "Monero"
Set-TextValue "C29" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D29" "157.03"
Set-TextValue "E29" "  -0.05%  "

Set-TextValue "B30" "BitcoinCash"
Set-TextValue "C30" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D30" "127.20"
Set-TextValue "E30" "  +0.07%  "

Set-TextValue "B31" "Stellar"
Set-TextValue "C31" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D31" "0.1058"
Set-TextValue "E31" "  -0.82%  "

Set-TextValue "B32" "ImmutableX"
Set-TextValue "C32" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D32" "1.064"
Set-TextValue "E32" "  -2.14%  "

Set-TextValue "B33" "Filecoin"
Set-TextValue "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "5.626"
Set-TextValue "E33" "  +0.16%  "

Set-TextValue "B34" "HuobiToken"
Set-TextValue "C34" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D34" "3.598"
Set-TextValue "E34" "  -0.52%  "

Set-TextValue "B35" "FraxShare"
Set-TextValue "C35" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D35" "9.605"
Set-TextValue "E35" "  -0.49%  "

Set-TextValue "B36" "Hedera"
Set-TextValue "C36" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D36" "0.06625"
Set-TextValue "E36" "  -0.20%  "

Set-TextValue "B37" "VeChain"
Set-TextValue "C37" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D37" "0.02413"
Set-TextValue "E37" "  -0.32%  "

Set-TextValue "B38" "Algorand"
Set-TextValue "C38" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D38" "0.2188"
Set-TextValue "E38" "  +0.65%  "

Set-TextValue "B39" "TrustWalletToken"
Set-TextValue "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D39" "1.286"
Set-TextValue "E39" "  +1.37%  "

Set-TextValue "B40" "ARBITRUM"
Set-TextValue "C40" "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D40" "1.212"
Set-TextValue "E40" "  -1.16%  "

Set-TextValue "B41" "TheSandbox"
Set-TextValue "C41" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D41" "0.6419"
Set-TextValue "E41" "  +0.47%  "

Set-TextValue "B42" "Aptos"
Set-TextValue "C42" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D42" "11.51"
Set-TextValue "E42" "  +0.93%  "

Set-TextValue "B43" "InternetComputer(DFINITY)"
Set-TextValue "C43" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D43" "4.931"
Set-TextValue "E43" "  -1.23%  "

Set-TextValue "B44" "Decentraland"
Set-TextValue "C44" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D44" "0.6052"
Set-TextValue "E44" "  +0.87%  "

Set-TextValue "B45" "EnergySwap"
Set-TextValue "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D45" "13.22"
Set-TextValue "E45" "  -0.64%  "

Set-TextValue "B46" "WEMIXTOKEN"
Set-TextValue "C46" "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue "D46" "1.276"
Set-TextValue "E46" "  +0.17%  "

Set-TextValue "D47" "3.668"
Set-TextValue "E47" "  -1.31%  "

Set-TextValue "B48" "EOS"
Set-TextValue "C48" "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue "D48" "1.245"
Set-TextValue "E48" "  +5.62%  "

Set-TextValue "B49" "NEARProtocol"
Set-TextValue "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D49" "2.000"
Set-TextValue "E49" "  -1.42%  "

Set-TextValue "B50" "Quant"
Set-TextValue "C50" "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue "D50" "121.36"
Set-TextValue "E50" "  -0.87%  "

Set-TextValue "B51" "Aave"
Set-TextValue "C51" "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D51" "79.63"
Set-TextValue "E51" "  +2.14%  "
